$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Alluvial for Mapping" (columns V = Max Cr, W = Max Date) ---
$ws1 = $wb.Worksheets.Item("Alluvial for Mapping")

$ws1.Range("V6").NumberFormat = "@"
$ws1.Range("V6").Value = "5.5"
$ws1.Range("W6").NumberFormat = "@"
$ws1.Range("W6").Value = "2008-02-11"
$ws1.Range("V7").NumberFormat = "@"
$ws1.Range("V7").Value = "No Detect Data"
$ws1.Range("W7").NumberFormat = "@"
$ws1.Range("W7").Value = "No Detect Data"
$ws1.Range("V13").NumberFormat = "@"
$ws1.Range("V13").Value = "6.69"
$ws1.Range("W13").NumberFormat = "@"
$ws1.Range("W13").Value = "2010-07-06"
$ws1.Range("V14").NumberFormat = "@"
$ws1.Range("V14").Value = "4.81"
$ws1.Range("W14").NumberFormat = "@"
$ws1.Range("W14").Value = "2010-07-07"
$ws1.Range("V15").NumberFormat = "@"
$ws1.Range("V15").Value = "No Detect Data"
$ws1.Range("W15").NumberFormat = "@"
$ws1.Range("W15").Value = "No Detect Data"
$ws1.Range("V16").NumberFormat = "@"
$ws1.Range("V16").Value = "4.69"
$ws1.Range("W16").NumberFormat = "@"
$ws1.Range("W16").Value = "2010-07-08"
$ws1.Range("V18").NumberFormat = "@"
$ws1.Range("V18").Value = "5.63"
$ws1.Range("W18").NumberFormat = "@"
$ws1.Range("W18").Value = "2010-07-09"
$ws1.Range("V19").NumberFormat = "@"
$ws1.Range("V19").Value = "2.2"
$ws1.Range("W19").NumberFormat = "@"
$ws1.Range("W19").Value = "2007-03-13"
$ws1.Range("V20").NumberFormat = "@"
$ws1.Range("V20").Value = "2.26"
$ws1.Range("W20").NumberFormat = "@"
$ws1.Range("W20").Value = "2013-06-05"
$ws1.Range("V21").NumberFormat = "@"
$ws1.Range("V21").Value = "No Detect Data"
$ws1.Range("W21").NumberFormat = "@"
$ws1.Range("W21").Value = "No Detect Data"
$ws1.Range("V22").NumberFormat = "@"
$ws1.Range("V22").Value = "5.7"
$ws1.Range("W22").NumberFormat = "@"
$ws1.Range("W22").Value = "2006-08-03"
$ws1.Range("V23").NumberFormat = "@"
$ws1.Range("V23").Value = "0.885"
$ws1.Range("W23").NumberFormat = "@"
$ws1.Range("W23").Value = "2003-07-09"
$ws1.Range("V26").NumberFormat = "@"
$ws1.Range("V26").Value = "4.9"
$ws1.Range("W26").NumberFormat = "@"
$ws1.Range("W26").Value = "2007-07-23"
$ws1.Range("V27").NumberFormat = "@"
$ws1.Range("V27").Value = "No Detect Data"
$ws1.Range("W27").NumberFormat = "@"
$ws1.Range("W27").Value = "No Detect Data"
$ws1.Range("V28").NumberFormat = "@"
$ws1.Range("V28").Value = "3.6"
$ws1.Range("W28").NumberFormat = "@"
$ws1.Range("W28").Value = "2007-07-19"
$ws1.Range("V29").NumberFormat = "@"
$ws1.Range("V29").Value = "1.6"
$ws1.Range("W29").NumberFormat = "@"
$ws1.Range("W29").Value = "2007-08-03"
$ws1.Range("V30").NumberFormat = "@"
$ws1.Range("V30").Value = "1.3"
$ws1.Range("W30").NumberFormat = "@"
$ws1.Range("W30").Value = "2005-05-10"
$ws1.Range("V31").NumberFormat = "@"
$ws1.Range("V31").Value = "1.02"
$ws1.Range("W31").NumberFormat = "@"
$ws1.Range("W31").Value = "2003-07-09"
$ws1.Range("V33").NumberFormat = "@"
$ws1.Range("V33").Value = "2.6"
$ws1.Range("W33").NumberFormat = "@"
$ws1.Range("W33").Value = "2007-08-02"
$ws1.Range("V34").NumberFormat = "@"
$ws1.Range("V34").Value = "6.72"
$ws1.Range("W34").NumberFormat = "@"
$ws1.Range("W34").Value = "2010-02-22"
$ws1.Range("V35").NumberFormat = "@"
$ws1.Range("V35").Value = "3.4"
$ws1.Range("W35").NumberFormat = "@"
$ws1.Range("W35").Value = "2008-06-17"
$ws1.Range("V36").NumberFormat = "@"
$ws1.Range("V36").Value = "3.8"
$ws1.Range("W36").NumberFormat = "@"
$ws1.Range("W36").Value = "2008-06-16"
$ws1.Range("V37").NumberFormat = "@"
$ws1.Range("V37").Value = "2.8"
$ws1.Range("W37").NumberFormat = "@"
$ws1.Range("W37").Value = "2006-12-11"
$ws1.Range("V38").NumberFormat = "@"
$ws1.Range("V38").Value = "5.56"
$ws1.Range("W38").NumberFormat = "@"
$ws1.Range("W38").Value = "2010-06-07"
$ws1.Range("V39").NumberFormat = "@"
$ws1.Range("V39").Value = "No Detect Data"
$ws1.Range("W39").NumberFormat = "@"
$ws1.Range("W39").Value = "No Detect Data"
$ws1.Range("V40").NumberFormat = "@"
$ws1.Range("V40").Value = "2.68"
$ws1.Range("V41").NumberFormat = "@"
$ws1.Range("V41").Value = "No Detect Data"
$ws1.Range("W41").NumberFormat = "@"
$ws1.Range("W41").Value = "No Detect Data"
$ws1.Range("V42").NumberFormat = "@"
$ws1.Range("V42").Value = "5.0"
$ws1.Range("W42").NumberFormat = "@"
$ws1.Range("W42").Value = "2008-12-18"
$ws1.Range("V43").NumberFormat = "@"
$ws1.Range("V43").Value = "2.74"
$ws1.Range("W43").NumberFormat = "@"
$ws1.Range("W43").Value = "2009-09-14"
$ws1.Range("V44").NumberFormat = "@"
$ws1.Range("V44").Value = "9.1"
$ws1.Range("W44").NumberFormat = "@"
$ws1.Range("W44").Value = "2008-09-17"
$ws1.Range("V45").NumberFormat = "@"
$ws1.Range("V45").Value = "2.0"
$ws1.Range("W45").NumberFormat = "@"
$ws1.Range("W45").Value = "2005-08-24"
$ws1.Range("V46").NumberFormat = "@"
$ws1.Range("V46").Value = "5.8"
$ws1.Range("W46").NumberFormat = "@"
$ws1.Range("W46").Value = "2007-06-25"
$ws1.Range("V47").NumberFormat = "@"
$ws1.Range("V47").Value = "3.6"
$ws1.Range("W47").NumberFormat = "@"
$ws1.Range("W47").Value = "2007-12-07"
$ws1.Range("V48").NumberFormat = "@"
$ws1.Range("V48").Value = "No Detect Data"
$ws1.Range("W48").NumberFormat = "@"
$ws1.Range("W48").Value = "No Detect Data"

# --- Sheet 2: "Alluvial Exhibit" (columns G = Max Cr [ug/L], H = Date of Max) ---
$ws2 = $wb.Worksheets.Item("Alluvial Exhibit")

$ws2.Range("G9").NumberFormat = "@"
$ws2.Range("G9").Value = "5.5"
$ws2.Range("H9").NumberFormat = "@"
$ws2.Range("H9").Value = "2/11/08"
$ws2.Range("G10").NumberFormat = "@"
$ws2.Range("G10").Value = "NA"
$ws2.Range("H10").NumberFormat = "@"
$ws2.Range("H10").Value = "NA"
$ws2.Range("G16").NumberFormat = "@"
$ws2.Range("G16").Value = "6.69"
$ws2.Range("H16").NumberFormat = "@"
$ws2.Range("H16").Value = "7/6/10"
$ws2.Range("G17").NumberFormat = "@"
$ws2.Range("G17").Value = "4.81"
$ws2.Range("H17").NumberFormat = "@"
$ws2.Range("H17").Value = "7/7/10"
$ws2.Range("G18").NumberFormat = "@"
$ws2.Range("G18").Value = "NA"
$ws2.Range("H18").NumberFormat = "@"
$ws2.Range("H18").Value = "NA"
$ws2.Range("G19").NumberFormat = "@"
$ws2.Range("G19").Value = "4.69"
$ws2.Range("H19").NumberFormat = "@"
$ws2.Range("H19").Value = "7/8/10"
$ws2.Range("G21").NumberFormat = "@"
$ws2.Range("G21").Value = "5.63"
$ws2.Range("H21").NumberFormat = "@"
$ws2.Range("H21").Value = "7/9/10"
$ws2.Range("G22").NumberFormat = "@"
$ws2.Range("G22").Value = "2.2"
$ws2.Range("H22").NumberFormat = "@"
$ws2.Range("H22").Value = "3/13/07"
$ws2.Range("G24").NumberFormat = "@"
$ws2.Range("G24").Value = "2.26"
$ws2.Range("H24").NumberFormat = "@"
$ws2.Range("H24").Value = "6/5/13"
$ws2.Range("G25").NumberFormat = "@"
$ws2.Range("G25").Value = "NA"
$ws2.Range("H25").NumberFormat = "@"
$ws2.Range("H25").Value = "NA"
$ws2.Range("G26").NumberFormat = "@"
$ws2.Range("G26").Value = "5.7"
$ws2.Range("H26").NumberFormat = "@"
$ws2.Range("H26").Value = "8/3/06"
$ws2.Range("G27").NumberFormat = "@"
$ws2.Range("G27").Value = "0.885"
$ws2.Range("H27").NumberFormat = "@"
$ws2.Range("H27").Value = "7/9/03"
$ws2.Range("G30").NumberFormat = "@"
$ws2.Range("G30").Value = "4.9"
$ws2.Range("H30").NumberFormat = "@"
$ws2.Range("H30").Value = "7/23/07"
$ws2.Range("G31").NumberFormat = "@"
$ws2.Range("G31").Value = "NA"
$ws2.Range("H31").NumberFormat = "@"
$ws2.Range("H31").Value = "NA"
$ws2.Range("G32").NumberFormat = "@"
$ws2.Range("G32").Value = "3.6"
$ws2.Range("H32").NumberFormat = "@"
$ws2.Range("H32").Value = "7/19/07"
$ws2.Range("G33").NumberFormat = "@"
$ws2.Range("G33").Value = "1.6"
$ws2.Range("H33").NumberFormat = "@"
$ws2.Range("H33").Value = "8/3/07"
$ws2.Range("G34").NumberFormat = "@"
$ws2.Range("G34").Value = "1.3"
$ws2.Range("H34").NumberFormat = "@"
$ws2.Range("H34").Value = "5/10/05"
$ws2.Range("G35").NumberFormat = "@"
$ws2.Range("G35").Value = "1.02"
$ws2.Range("H35").NumberFormat = "@"
$ws2.Range("H35").Value = "7/9/03"
$ws2.Range("G37").NumberFormat = "@"
$ws2.Range("G37").Value = "2.6"
$ws2.Range("H37").NumberFormat = "@"
$ws2.Range("H37").Value = "8/2/07"
$ws2.Range("G38").NumberFormat = "@"
$ws2.Range("G38").Value = "6.72"
$ws2.Range("H38").NumberFormat = "@"
$ws2.Range("H38").Value = "2/22/10"
$ws2.Range("G39").NumberFormat = "@"
$ws2.Range("G39").Value = "3.4"
$ws2.Range("H39").NumberFormat = "@"
$ws2.Range("H39").Value = "6/17/08"
$ws2.Range("G40").NumberFormat = "@"
$ws2.Range("G40").Value = "3.8"
$ws2.Range("H40").NumberFormat = "@"
$ws2.Range("H40").Value = "6/16/08"
$ws2.Range("G41").NumberFormat = "@"
$ws2.Range("G41").Value = "2.8"
$ws2.Range("H41").NumberFormat = "@"
$ws2.Range("H41").Value = "12/11/06"
$ws2.Range("G42").NumberFormat = "@"
$ws2.Range("G42").Value = "5.56"
$ws2.Range("H42").NumberFormat = "@"
$ws2.Range("H42").Value = "6/7/10"
$ws2.Range("G43").NumberFormat = "@"
$ws2.Range("G43").Value = "NA"
$ws2.Range("H43").NumberFormat = "@"
$ws2.Range("H43").Value = "NA"
$ws2.Range("G44").NumberFormat = "@"
$ws2.Range("G44").Value = "2.68"
$ws2.Range("G45").NumberFormat = "@"
$ws2.Range("G45").Value = "NA"
$ws2.Range("H45").NumberFormat = "@"
$ws2.Range("H45").Value = "NA"
$ws2.Range("G46").NumberFormat = "@"
$ws2.Range("G46").Value = "5"
$ws2.Range("H46").NumberFormat = "@"
$ws2.Range("H46").Value = "12/18/08"
$ws2.Range("G47").NumberFormat = "@"
$ws2.Range("G47").Value = "2.74"
$ws2.Range("H47").NumberFormat = "@"
$ws2.Range("H47").Value = "9/14/09"
$ws2.Range("G48").NumberFormat = "@"
$ws2.Range("G48").Value = "9.1"
$ws2.Range("H48").NumberFormat = "@"
$ws2.Range("H48").Value = "9/17/08"
$ws2.Range("G49").NumberFormat = "@"
$ws2.Range("G49").Value = "2"
$ws2.Range("H49").NumberFormat = "@"
$ws2.Range("H49").Value = "8/24/05"
$ws2.Range("G50").NumberFormat = "@"
$ws2.Range("G50").Value = "5.8"
$ws2.Range("H50").NumberFormat = "@"
$ws2.Range("H50").Value = "6/25/07"
$ws2.Range("G51").NumberFormat = "@"
$ws2.Range("G51").Value = "3.6"
$ws2.Range("H51").NumberFormat = "@"
$ws2.Range("H51").Value = "12/7/07"
$ws2.Range("G52").NumberFormat = "@"
$ws2.Range("G52").Value = "NA"
$ws2.Range("H52").NumberFormat = "@"
$ws2.Range("H52").Value = "NA"
